# Apply the changes described by the diff to Tab_5a_Indikatoren.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (columns D, E, F, I -> indexes 4, 5, 6, 9) ---
# Target OOXML <col> widths (stored in 1/256 character units) are:
#   col 4 -> 15.37890625
#   col 5 -> 18.6015625
#   col 6 -> 17.28515625
#   col 9 -> 24.31640625
# The ColumnWidth COM property is rounded to the nearest pixel on a 7px
# "Maximum Digit Width" grid, so we pick the ColumnWidth value whose
# rounded pixel width is as close as possible to each target.
$ws.Columns.Item(4).ColumnWidth = 14.714285714285714
$ws.Columns.Item(5).ColumnWidth = 17.857142857142858
$ws.Columns.Item(6).ColumnWidth = 16.571428571428573
$ws.Columns.Item(9).ColumnWidth = 23.571428571428573

# --- Row 5 (Z02_B01_P01_Ib02_I01 / Organic farming): 20% -> 30% by 2030 ---
$ws.Range("H5").Value = "Erhöhung des Anteils des ökologischen Landbaus an der landwirtschaftlich genutzten Fläche auf 30 % bis 2030"
$ws.Range("I5").Value = "Increase the proportion of organically farmed agricultural land to 30% by 2030"
$ws.Range("J5").Value = "Erhöhung des Anteils auf 30 % bis 2030"
$ws.Range("K5").Value = "Increase the proportion to 30% by 2030"

# --- Row 19 (All-day care for 3-5 year-olds): drop 2020 milestone, keep 70% by 2030 ---
$ws.Range("H19").Value = "Anstieg auf 70 % bis 2030"
$ws.Range("I19").Value = "Increase to 70% by 2030"
$ws.Range("J19").Value = "Anstieg auf 70 % bis 2030"
$ws.Range("K19").Value = "Increase to 70% by 2030"

# --- Row 33 (Share of electricity from renewables): updated to at least 80% by 2030 ---
$ws.Range("H33").Value = "Anstieg auf mindestens 80 % bis 2030"
$ws.Range("I33").Value = "Increase to at least 80% by 2030"
$ws.Range("J33").Value = "Anstieg auf 80 % bis 2030"
$ws.Range("K33").Value = "Increase to 80% by 2030"

# --- Row 62 (Greenhouse gas emissions): updated targets/years ---
$ws.Range("H62").Value = "Minderung um mindestens 65 % bis 2030, um mindestens 88 % bis 2040; Erreichung der Treibhausgasneutralität bis 2045"
$ws.Range("I62").Value = "Reduce by at least 65% by 2030 and by at least 88% by 2040; greenhouse gas neutrality to be achieved by 2045"
$ws.Range("J62").Value = "Minderung um 65 % bis 2030"
$ws.Range("K62").Value = "Reduce by 65% by 2030"

# --- Row 63 (Climate finance payments): new financing target wording ---
$ws.Range("H63").Value = "Erhöhung der internationalen Klimafinanzierung auf mindestens 6 Mrd. Euro bis spätestens 2025."
$ws.Range("I63").Value = "Increase international climate finance to at least 6 billion euros by 2025 at the latest."
$ws.Range("J63").Value = "Erhöhung auf 6 Mrd. Euro bis 2025"
$ws.Range("K63").Value = "Increase to 6 Mrd. euro by 2025"
